# Insert a new data row at row 641 (shifting the existing rows 641-693
# down to 642-694) and populate the new row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 641 (and everything below it) down by one row.
$ws.Rows.Item(641).Insert()

# Populate the newly inserted row 641 with the new observation.
$ws.Cells.Item(641, 1).Value  = 3
$ws.Cells.Item(641, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(641, 3).Value  = "Coquimbo"
$ws.Cells.Item(641, 4).Value2 = 45223
$ws.Cells.Item(641, 5).Value  = 5
$ws.Cells.Item(641, 6).Value  = 100112031
$ws.Cells.Item(641, 7).Value  = "Poroto verde"
$ws.Cells.Item(641, 8).Value  = "Sin especificar"
$ws.Cells.Item(641, 9).Value  = "Primera"
$ws.Cells.Item(641, 10).Value = 60
$ws.Cells.Item(641, 11).Value = 31000
$ws.Cells.Item(641, 12).Value = 32000
$ws.Cells.Item(641, 13).Value = 31500
$ws.Cells.Item(641, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(641, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(641, 16).Value = 1260
$ws.Cells.Item(641, 17).Value = 25
$ws.Cells.Item(641, 18).Value = "Hortaliza"

# Make sure the D641 cell keeps the date number format used by the rest
# of column D (same style as D640/D642, numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Cells.Item(641, 4).NumberFormat = $ws.Cells.Item(642, 4).NumberFormat
